# Repo: Repository.TestGitCommitChangesOutsideWebstudio / Main.xlsx
# "update file with jgit" - the Rules sheet's R10 greeting cell (E8) was
# changed from "Good Morning" to "GIT UPDATE", and the sheet's selection
# ended up on that same cell (E8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

# Leave the active selection on E8 (matches the <selection> added to the
# sheetView in the target workbook).
$null = $ws.Range("E8").Select()
